# Essential changes to datasets
# Prefix the Supplier names (column A) with a letter code (A-E) for each
# distinct supplier group, while keeping the rest of the data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$prefixes = @{
    "Mono Packaging Materials" = "A Mono Packaging Materials"
    "Trio PET PLC"             = "B Trio PET PLC"
    "Miami Oranges"            = "C Miami Oranges"
    "NO8DO Mango"              = "D NO8DO Mango"
    "Seitan Vitamins"          = "E Seitan Vitamins"
}

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Text
    if ($prefixes.ContainsKey($current)) {
        $cell.Value = $prefixes[$current]
    }
}

# Update the selection to match the edited workbook state
$ws.Range("A2:A21").Select()
